$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 13.02000045776367
$ws.Range("F2").Value = 15.73999977111816
$ws.Range("G2").Value = 11.94999980926514
$ws.Range("H2").Value = 274203783
$ws.Range("I2").Value = "PONY"

# Row 3
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 13.02000045776367
$ws.Range("F3").Value = 15.73999977111816
$ws.Range("G3").Value = 11.94999980926514
$ws.Range("H3").Value = 274203783
$ws.Range("I3").Value = "PONY"

# Row 4
$ws.Range("D4").Value = 15
$ws.Range("E4").Value = 13.02000045776367
$ws.Range("F4").Value = 15.73999977111816
$ws.Range("G4").Value = 11.94999980926514
$ws.Range("H4").Value = 274203783
$ws.Range("I4").Value = "PONY"

# Row 5
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 13.02000045776367
$ws.Range("F5").Value = 15.73999977111816
$ws.Range("G5").Value = 11.94999980926514
$ws.Range("H5").Value = 274203783
$ws.Range("I5").Value = "PONY"

# Row 6
$ws.Range("D6").Value = 15.34000015258789
$ws.Range("E6").Value = 13.15999984741211
$ws.Range("F6").Value = 16.25
$ws.Range("G6").Value = 12.51000022888184
$ws.Range("H6").Value = 274203783
$ws.Range("I6").Value = "PONY"

# Row 7
$ws.Range("D7").Value = 8.840000152587891
$ws.Range("E7").Value = 9.090000152587891
$ws.Range("F7").Value = 11.34500026702881
$ws.Range("G7").Value = 4.110000133514404
$ws.Range("H7").Value = 274203783
$ws.Range("I7").Value = "PONY"

# Row 8
$ws.Range("D8").Value = 12.69999980926514
$ws.Range("E8").Value = 13.4399995803833
$ws.Range("F8").Value = 16.38999938964844
$ws.Range("G8").Value = 11.77000045776367
$ws.Range("H8").Value = 274203783
$ws.Range("I8").Value = "PONY"
